# Auto-generated edit script: updates leve profit calculation cells
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) to reflect
# refreshed market-board pricing data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 15214.429
$ws.Range("I12").Value = 33433.668
$ws.Range("J12").Value = 1550
$ws.Range("K12").Value = 33433.668
$ws.Range("L12").Value = 1550
$ws.Range("M12").Value = -33263.668
$ws.Range("N12").Value = -1890

$ws.Range("H113").Value = 3288.2307
$ws.Range("I113").Value = 2339.8
$ws.Range("J113").Value = 3881
$ws.Range("K113").Value = 2339.8
$ws.Range("L113").Value = 3881
$ws.Range("M113").Value = 914.1999999999998
$ws.Range("N113").Value = -10389

$ws.Range("H137").Value = 733899.5
$ws.Range("I137").Value = 1927.8572
$ws.Range("J137").Value = 1638099.8
$ws.Range("K137").Value = 5783.571599999999
$ws.Range("L137").Value = 4914299.4
$ws.Range("M137").Value = -3233.571599999999
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 2698.0632
$ws.Range("I138").Value = 1113.8695
$ws.Range("J138").Value = 4185.265
$ws.Range("K138").Value = 3341.6085
$ws.Range("L138").Value = 12555.795
$ws.Range("M138").Value = 1798.3915
$ws.Range("N138").Value = -22835.795

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1127.5834
$ws.Range("I2").Value = 1185.0416
$ws.Range("J2").Value = 1012.6667
$ws.Range("K2").Value = 1185.0416
$ws.Range("L2").Value = 1012.6667
$ws.Range("M2").Value = -1072.0416
$ws.Range("N2").Value = -1238.6667

$ws.Range("H32").Value = 7673.7114
$ws.Range("I32").Value = 6149.638
$ws.Range("J32").Value = 22000
$ws.Range("K32").Value = 6149.638
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = -5862.638
$ws.Range("N32").ClearContents()

$ws.Range("H61").Value = 5454.5454
$ws.Range("I61").Value = 3720.9644
$ws.Range("J61").Value = 15162.6
$ws.Range("K61").Value = 3720.9644
$ws.Range("L61").Value = 15162.6
$ws.Range("M61").Value = -3508.9644
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 6312.0347
$ws.Range("I74").Value = 3305.625
$ws.Range("J74").Value = 20742.8
$ws.Range("K74").Value = 3305.625
$ws.Range("L74").Value = 20742.8
$ws.Range("M74").Value = -2431.625
$ws.Range("N74").Value = -22490.8

$ws.Range("H77").Value = 6312.0347
$ws.Range("I77").Value = 3305.625
$ws.Range("J77").Value = 20742.8
$ws.Range("K77").Value = 16528.125
$ws.Range("L77").Value = 103714
$ws.Range("M77").Value = -12160.125
$ws.Range("N77").Value = -112450

$ws.Range("H116").Value = 1127.5834
$ws.Range("I116").Value = 1185.0416
$ws.Range("J116").Value = 1012.6667
$ws.Range("K116").Value = 1185.0416
$ws.Range("L116").Value = 1012.6667
$ws.Range("M116").Value = 1108.9584
$ws.Range("N116").Value = -5600.6667

$ws.Range("H132").Value = 7567.6523
$ws.Range("I132").Value = 1977.091
$ws.Range("J132").Value = 12692.333
$ws.Range("K132").Value = 5931.272999999999
$ws.Range("L132").Value = 38076.999
$ws.Range("M132").Value = -3401.272999999999
$ws.Range("N132").Value = -43136.999

$ws.Range("H136").Value = 5454.5454
$ws.Range("I136").Value = 3720.9644
$ws.Range("J136").Value = 15162.6
$ws.Range("K136").Value = 11162.8932
$ws.Range("L136").Value = 45487.8
$ws.Range("M136").Value = -8612.893199999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1127.5834
$ws.Range("I3").Value = 1185.0416
$ws.Range("J3").Value = 1012.6667
$ws.Range("K3").Value = 1185.0416
$ws.Range("L3").Value = 1012.6667
$ws.Range("M3").Value = -1071.0416
$ws.Range("N3").Value = -1240.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2458.8333
$ws.Range("I31").Value = 1797
$ws.Range("J31").Value = 7753.5
$ws.Range("K31").Value = 1797
$ws.Range("L31").Value = 7753.5
$ws.Range("M31").Value = -1502
$ws.Range("N31").Value = -8343.5

$ws.Range("H34").Value = 2458.8333
$ws.Range("I34").Value = 1797
$ws.Range("J34").Value = 7753.5
$ws.Range("K34").Value = 1797
$ws.Range("L34").Value = 7753.5
$ws.Range("M34").Value = -1595
$ws.Range("N34").Value = -8157.5

$ws.Range("H86").Value = 1798.5
$ws.Range("I86").Value = 990
$ws.Range("J86").Value = 2068
$ws.Range("K86").Value = 990
$ws.Range("L86").Value = 2068
$ws.Range("M86").Value = 133
$ws.Range("N86").Value = -4314

$ws.Range("H89").Value = 1798.5
$ws.Range("I89").Value = 990
$ws.Range("J89").Value = 2068
$ws.Range("K89").Value = 4950
$ws.Range("L89").Value = 10340
$ws.Range("M89").Value = 666
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1939
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3696
$ws.Range("N77").ClearContents()

$ws.Range("H96").Value = 6000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 18000
$ws.Range("N96").Value = -22118
$ws.Range("M96").ClearContents()

$ws.Range("H131").Value = 1097.5135
$ws.Range("I131").Value = 1413.2142
$ws.Range("J131").Value = 905.34784
$ws.Range("K131").Value = 4239.642599999999
$ws.Range("L131").Value = 2716.04352
$ws.Range("M131").Value = 800.3574000000008
$ws.Range("N131").Value = -12796.04352

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 3951.6667
$ws.Range("I12").Value = 403
$ws.Range("J12").Value = 4661.4
$ws.Range("K12").Value = 403
$ws.Range("L12").Value = 4661.4
$ws.Range("M12").Value = -233
$ws.Range("N12").Value = -5001.4

$ws.Range("H22").Value = 718.0909
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 738.1177
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 738.1177
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1328.1177

$ws.Range("H27").Value = 718.0909
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 738.1177
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 738.1177
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -952.1177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 375
$ws.Range("I126").Value = 322.5
$ws.Range("J126").Value = 1005
$ws.Range("K126").Value = 967.5
$ws.Range("L126").Value = 3015
$ws.Range("M126").Value = 1502.5
$ws.Range("N126").Value = -7955

$ws.Range("H132").Value = 1310.425
$ws.Range("I132").Value = 1530.9584
$ws.Range("J132").Value = 979.625
$ws.Range("K132").Value = 4592.8752
$ws.Range("L132").Value = 2938.875
$ws.Range("M132").Value = -2062.8752
$ws.Range("N132").Value = -7998.875

$ws.Range("H136").Value = 4880.8
$ws.Range("I136").Value = 2515
$ws.Range("J136").Value = 7064.615
$ws.Range("K136").Value = 7545
$ws.Range("L136").Value = 21193.845
$ws.Range("M136").Value = -4995
$ws.Range("N136").Value = -26293.845
